# Anna (tasks, problems, how solving problems)
#
# Fix a typo in the ConnectionPool write-up (D12), replace the
# "Special problems have arisen" note (C15) with the sequence-diagram
# follow-up question, add the matching "how we solved it" note (D15),
# and move the sheet selection from G15 to F14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "The fact that our ConnectionPool  is safe-thread confirmed using singleton pattern and the fact it is multithreaded singleton  confirmed pool setting . Find information   about singletone and ThreadLocal class."

$ws.Range("D15").Value = "Find tutorials about it."
$ws.Range("C15").Value = "How sequence diagram to depict a link to another sequence diagram."

$ws.Range("F14").Select()
